# Update the "想去人数" (want-to-go count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 127
$wsExhibit.Range("F4").Value = 676

# Sheet "全部类型" (All types) - same two events appear here on rows 4 and 5
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 127
$wsAll.Range("F5").Value = 676
